# Update gh-pages output-generated numbers ("想去人数" counts) across the
# 4 worksheets, and drop the now-removed "本地生活" row for the expired
# "贰伊Lolita茶会" event.
#
# Sheet order (tab order == workbook.xml <sheets>):
#   1 = 展览 (exhibitions)
#   2 = 演出 (performances)
#   3 = 本地生活 (local life)
#   4 = 全部类型 (all types, combined)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 -----------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$sheet1Updates = @{
    2  = 16
    5  = 5051
    6  = 5051
    7  = 85
    9  = 498
    11 = 1143
    12 = 695
    13 = 4869
    14 = 21
    15 = 52
    16 = 68
    18 = 209
    19 = 94
    21 = 3732
    24 = 3583
    28 = 192
    31 = 103
    32 = 104
    36 = 6315
    37 = 998
    38 = 475
    39 = 94
    42 = 1295
    43 = 150
    44 = 623
    46 = 2179
    49 = 753
    50 = 896
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet 2: 演出 -------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$sheet2Updates = @{
    2  = 14
    9  = 72
    23 = 794
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("F$row").Value = $sheet2Updates[$row]
}

# --- Sheet 3: 本地生活 ---------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 215
# The "2024-06-16 北京·贰伊Lolita茶会" event (row 3) dropped out of this feed;
# delete the row entirely (dimension auto-shrinks from A1:I3 to A1:I2).
$ws3.Rows.Item(3).Delete()

# --- Sheet 4: 全部类型 ---------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$sheet4Updates = @{
    2  = 215
    3  = 14
    5  = 16
    10 = 5051
    11 = 5051
    12 = 85
    15 = 72
    16 = 695
    17 = 4869
    18 = 21
    19 = 52
    20 = 68
    22 = 94
    24 = 3583
    27 = 192
    30 = 103
    31 = 104
    36 = 6315
    37 = 998
    38 = 94
    40 = 1295
    41 = 150
    42 = 623
    44 = 2179
    48 = 753
    49 = 896
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
